$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the report-month labels: "January 2017" -> "February 2017" ---
$ws.Range("A2").Value   = "Short-Term Energy Outlook, February 2017"
$ws.Range("A100").Value = "Source: Short-Term Energy Outlook, February 2017."

# --- Refresh the monthly price history / forecast (columns B & C, rows 74-99) ---
$priceData = @(
    @(10.76,               2.6295359999999999),
    @(9.8132370000000009,  3.7059120000000001),
    @(10.1983,             3.4097279999999999),
    @(9.8720669999999995,  3.505128),
    @(10.225770000000001,  3.5003669999999998),
    @(11.111269999999999,  3.5006550000000001),
    @(13.0837,             3.4730349999999999),
    @(15.398870000000001,  3.5011049999999999),
    @(16.615259999999999,  3.527911),
    @(17.446429999999999,  3.5169299999999999),
    @(16.494810000000001,  3.5100419999999999),
    @(13.49583,            3.5418059999999998),
    @(11.117000000000001,  3.6551719999999999),
    @(10.23742,            3.81962),
    @(10.048310000000001,  3.888674),
    @(10.112299999999999,  3.9194100000000001),
    @(10.440709999999999,  3.8827569999999998),
    @(11.35529,            3.8517160000000001),
    @(13.36858,            3.7765870000000001),
    @(15.692629999999999,  3.7726220000000001),
    @(16.923449999999999,  3.7675010000000002),
    @(17.79383,            3.726788),
    @(16.817240000000002,  3.6974629999999999),
    @(13.798019999999999,  3.71143),
    @(11.380039999999999,  3.8117909999999999),
    @(10.48146,            3.967476)
)

$startRow = 74
for ($i = 0; $i -lt $priceData.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $priceData[$i][0]
    $ws.Cells.Item($row, 3).Value = $priceData[$i][1]
}

# --- Forecast-marker helper column: shift the "start of forecast" index 48 -> 49 ---
$ws.Range("A103").Value = 49
$ws.Range("A104").Value = 49
